$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 9478
$ws.Range("C2").Value = 9456
$ws.Range("D2").Value = 9456
$ws.Range("F2").Value = 0.9976788351972991
$ws.Range("G2").Value = 0.09258994317483452
$ws.Range("H2").Value = 0.092375026657653
$ws.Range("I2").Value = 47112440.0389384
$ws.Range("J2").Value = 17356555.7788052
$ws.Range("L2").Value = 17356555.7788052
$ws.Range("M2").Value = 64468995.81774361
$ws.Range("N2").Value = 813109489.4172001
$ws.Range("O2").Value = 794833158.7732
$ws.Range("P2").Value = 0.02134590237194943
$ws.Range("Q2").Value = 0.02183672835893573

# Row 3
$ws.Range("B3").Value = 9762
$ws.Range("C3").Value = 9738
$ws.Range("D3").Value = 9738
$ws.Range("F3").Value = 0.9975414874001229
$ws.Range("G3").Value = 0.09761071208907925
$ws.Range("H3").Value = 0.09737073492352527
$ws.Range("I3").Value = 55505811.34755692
$ws.Range("J3").Value = 21282186.73698258
$ws.Range("L3").Value = 21282186.73698258
$ws.Range("M3").Value = 76787998.0845395
$ws.Range("N3").Value = 863794315.9995871
$ws.Range("O3").Value = 845379885.828617
$ws.Range("P3").Value = 0.0246380259082334
$ws.Range("Q3").Value = 0.0251747020407546

# Row 4
$ws.Range("B4").Value = 10046
$ws.Range("C4").Value = 10024
$ws.Range("D4").Value = 10024
$ws.Range("F4").Value = 0.9978100736611587
$ws.Range("G4").Value = 0.1008389917673471
$ws.Range("H4").Value = 0.1006181618032936
$ws.Range("I4").Value = 62994171.20052955
$ws.Range("J4").Value = 24707287.94181236
$ws.Range("L4").Value = 24707287.94181236
$ws.Range("M4").Value = 87701459.1423419
$ws.Range("N4").Value = 916693029.3738154
$ws.Range("O4").Value = 897970254.2003546
$ws.Range("P4").Value = 0.02695262988820771
$ws.Range("Q4").Value = 0.02751459508401454

# Row 5
$ws.Range("B5").Value = 10337
$ws.Range("C5").Value = 10309
$ws.Range("D5").Value = 10309
$ws.Range("F5").Value = 0.9972912837380284
$ws.Range("G5").Value = 0.1007746111987947
$ws.Range("H5").Value = 0.1005016413706467
$ws.Range("I5").Value = 67918402.13472392
$ws.Range("J5").Value = 26836175.90669475
$ws.Range("L5").Value = 26836175.90669475
$ws.Range("M5").Value = 94754578.04141869
$ws.Range("N5").Value = 969816544.9170408
$ws.Range("O5").Value = 951407665.5136355
$ws.Range("P5").Value = 0.02767139419032117
$ws.Range("Q5").Value = 0.02820681068635992

# Row 6
$ws.Range("B6").Value = 10656
$ws.Range("C6").Value = 10631
$ws.Range("D6").Value = 10631
$ws.Range("F6").Value = 0.9976539039039038
$ws.Range("G6").Value = 0.1002504979019584
$ws.Range("H6").Value = 0.1000153006001989
$ws.Range("I6").Value = 73124236.81502223
$ws.Range("J6").Value = 29067604.94581254
$ws.Range("L6").Value = 29067604.94581254
$ws.Range("M6").Value = 102191841.7608348
$ws.Range("N6").Value = 1033114253.054143
$ws.Range("O6").Value = 1014594770.022316
$ws.Range("P6").Value = 0.02813590545274296
$ws.Range("Q6").Value = 0.02864947248365294
